# Update cryptos list with latest prices/volume figures.
#
# NOTE: the "Price" column stores plain-text strings (e.g. "27.902.08" or
# "0.9979"), not real numbers. Whenever the new price text would otherwise
# parse as a number (a single decimal point, no thousands separators),
# it's prefixed with a leading apostrophe - exactly like a user typing it
# in the Excel UI - so it is stored as text (preserving trailing zeros,
# leading zeros, etc.) instead of being auto-coerced into a numeric value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.931.34'
$ws.Range("E2").Value = '  -0.28%  '
$ws.Range("D3").Value = '1.911.67'
$ws.Range("E3").Value = '  +0.25%  '
$ws.Range("D4").Value = '''0.9987'
$ws.Range("E4").Value = '  -0.42%  '
$ws.Range("D5").Value = '''312.76'
$ws.Range("E5").Value = '  -1.70%  '
$ws.Range("D6").Value = '''0.9984'
$ws.Range("E6").Value = '  -0.48%  '
$ws.Range("D7").Value = '''0.5005'
$ws.Range("E7").Value = '  +3.78%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").Value = '''0.07290'
$ws.Range("E9").Value = '  -1.00%  '
$ws.Range("D10").Value = '''21.34'
$ws.Range("E10").Value = '  +2.51%  '
$ws.Range("D11").Value = '''0.9097'
$ws.Range("E11").Value = '  -2.44%  '
$ws.Range("D12").Value = '''0.07645'
$ws.Range("E12").Value = '  -1.36%  '
$ws.Range("D13").Value = '1.925.18'
$ws.Range("E13").Value = '  +0.95%  '
$ws.Range("D14").Value = '''5.476'
$ws.Range("E14").Value = '  -0.26%  '
$ws.Range("D15").Value = '''92.71'
$ws.Range("E15").Value = '  +1.03%  '
$ws.Range("D16").Value = '''0.9995'
$ws.Range("E16").Value = '  -0.49%  '
$ws.Range("D17").Value = '''0.000008735'
$ws.Range("E17").Value = '  -1.81%  '
$ws.Range("E18").Value = '  -0.46%  '
$ws.Range("D19").Value = '27.955.43'
$ws.Range("E19").Value = '  -0.31%  '
$ws.Range("D20").Value = '''14.67'
$ws.Range("E20").Value = '  -0.34%  '
$ws.Range("D21").Value = '''5.175'
$ws.Range("E21").Value = '  +0.62%  '
$ws.Range("D22").Value = '2.119.17'
$ws.Range("E22").Value = '  -1.20%  '
$ws.Range("D23").Value = '''10.88'
$ws.Range("D24").Value = '''6.603'
$ws.Range("E24").Value = '  -0.58%  '
$ws.Range("D25").Value = '''152.88'
$ws.Range("E25").Value = '  -2.45%  '
$ws.Range("E26").Value = '  -2.39%  '
$ws.Range("D27").Value = '''2.231'
$ws.Range("E27").Value = '  +5.27%  '
$ws.Range("D28").Value = '''18.40'
$ws.Range("E28").Value = '  -0.67%  '
$ws.Range("D29").Value = '''115.04'
$ws.Range("E29").Value = '  -2.12%  '
$ws.Range("D30").Value = '''4.906'
$ws.Range("E30").Value = '  -1.50%  '
$ws.Range("D31").Value = '''0.08989'
$ws.Range("E31").Value = '  +0.43%  '
$ws.Range("D32").Value = '''3.199'
$ws.Range("E32").Value = '  -1.42%  '
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").Value = '''0.7886'
$ws.Range("E33").Value = '  +2.13%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Value = '''4.815'
$ws.Range("E34").Value = '  +3.13%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '''1.233'
$ws.Range("E35").Value = '  -1.83%  '
$ws.Range("D36").Value = '''2.648'
$ws.Range("E36").Value = '  +1.47%  '
$ws.Range("D37").Value = '''0.02086'
$ws.Range("E37").Value = '  +1.77%  '
$ws.Range("D38").Value = '''3.061'
$ws.Range("E38").Value = '  +2.18%  '
$ws.Range("E39").Value = '  -1.43%  '
$ws.Range("D40").Value = '''0.5540'
$ws.Range("E40").Value = '  +0.38%  '
$ws.Range("D41").Value = '''0.05287'
$ws.Range("E41").Value = '  -0.12%  '
$ws.Range("E42").Value = '  -2.74%  '
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").Value = '''8.537'
$ws.Range("E43").Value = '  +0.41%  '
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = '''113.29'
$ws.Range("E44").Value = '  +2.60%  '
$ws.Range("D45").Value = '''0.1519'
$ws.Range("E45").Value = '  -0.46%  '
$ws.Range("E46").Value = '  -0.47%  '
$ws.Range("D47").Value = '''0.4832'
$ws.Range("E47").Value = '  -0.02%  '
$ws.Range("E48").Value = '  -0.49%  '
$ws.Range("E49").Value = '  -0.47%  '
$ws.Range("D50").Value = '''67.39'
$ws.Range("E50").Value = '  -0.82%  '
$ws.Range("D51").Value = '''0.06043'
$ws.Range("E51").Value = '  -0.62%  '
